$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Headers (M2:P2) mirror A2:D2 (Time / Control / 1 Day / 1 Week) ---
$ws.Range("A2:D2").Copy()
$ws.Range("M2").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("M2").Value2 = $ws.Range("A2").Value2
$ws.Range("N2").Value2 = $ws.Range("B2").Value2
$ws.Range("O2").Value2 = $ws.Range("C2").Value2
$ws.Range("P2").Value2 = $ws.Range("D2").Value2

# --- Row labels M3:M12 mirror A3:A12 ---
for ($r = 3; $r -le 12; $r++) {
    $ws.Range("A$r").Copy()
    $ws.Range("M$r").PasteSpecial(-4122)  # xlPasteFormats
    $ws.Range("M$r").Value2 = $ws.Range("A$r").Value2
}

# --- Percentage-difference columns N, O, P (rows 3-12) ---
# Style: same border/font as the data columns (B:D, H:J) but with a 0% number format.
for ($r = 3; $r -le 12; $r++) {
    foreach ($col in @("N", "O", "P")) {
        $dst = $ws.Range("$col$r")
        $ws.Range("B$r").Copy()
        $dst.PasteSpecial(-4122)  # xlPasteFormats
        $dst.NumberFormat = "0%"
    }
}

# N: ABS((B-H)/B) ; O: ABS((C-I)/C) ; P: ABS((D-J)/D)
# Where the divisor is zero Excel raises #DIV/0!; those cells were hand-replaced
# with the literal marker "X" rather than left as an error.
$formulas = @{
    3  = @($true, $true, $true)
    4  = @($true, $true, $true)
    5  = @($true, $true, $true)
    6  = @($true, $true, $true)
    7  = @($true, $true, $false)
    8  = @($true, $true, $false)
    9  = @($true, $true, $false)
    10 = @($false, $true, $false)
    11 = @($true, $true, $false)
    12 = @($true, $true, $false)
}

foreach ($r in 3..12) {
    $flags = $formulas[$r]

    if ($flags[0]) {
        $ws.Range("N$r").Formula = "=ABS((B$r-H$r)/B$r)"
    } else {
        $ws.Range("N$r").Value2 = "X"
    }

    if ($flags[1]) {
        $ws.Range("O$r").Formula = "=ABS((C$r-I$r)/C$r)"
    } else {
        $ws.Range("O$r").Value2 = "X"
    }

    if ($flags[2]) {
        $ws.Range("P$r").Formula = "=ABS((D$r-J$r)/D$r)"
    } else {
        $ws.Range("P$r").Value2 = "X"
    }
}

# --- View state: scrolled right a bit, selection on N17 ---
$ws.Range("N17").Select()
